$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.031.24"
$ws.Range("E2").Value = "  -3.71%  "
$ws.Range("D3").Value = "1.650.21"
$ws.Range("E3").Value = "  -5.37%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.9994"
$c.NumberFormat = "General"
$ws.Range("E4").Value = "  -0.04%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "236.98"
$c.NumberFormat = "General"
$ws.Range("E5").Value = "  -5.50%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.NumberFormat = "General"
$ws.Range("E6").Value = "  +0.02%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4829"
$c.NumberFormat = "General"
$ws.Range("E7").Value = "  -6.31%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.2621"
$c.NumberFormat = "General"
$ws.Range("E8").Value = "  -5.10%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.06027"
$c.NumberFormat = "General"
$ws.Range("E9").Value = "  -2.75%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.07198"
$c.NumberFormat = "General"
$ws.Range("E10").Value = "  -0.56%  "
$ws.Range("D11").Value = "1.650.47"
$ws.Range("E11").Value = "  -5.31%  "
$ws.Range("E12").Value = "  -2.68%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.6209"
$c.NumberFormat = "General"
$ws.Range("E13").Value = "  -4.82%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "4.580"
$c.NumberFormat = "General"
$ws.Range("E14").Value = "  -1.20%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "73.01"
$c.NumberFormat = "General"
$ws.Range("E15").Value = "  -6.31%  "
$ws.Range("E16").Value = "  +0.08%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.9988"
$c.NumberFormat = "General"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").Value = "25.017.68"
$ws.Range("E18").Value = "  -3.85%  "
$ws.Range("E19").Value = "  -3.03%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.000006630"
$c.NumberFormat = "General"
$ws.Range("E20").Value = "  -2.72%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "4.551"
$c.NumberFormat = "General"
$ws.Range("E21").Value = "  +5.71%  "
$ws.Range("D22").Value = "1.857.28"
$ws.Range("E22").Value = "  -5.57%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "8.619"
$c.NumberFormat = "General"
$ws.Range("E23").Value = "  -0.74%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "5.306"
$c.NumberFormat = "General"
$ws.Range("E24").Value = "  -1.62%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "132.07"
$c.NumberFormat = "General"
$ws.Range("E25").Value = "  -2.96%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "14.94"
$c.NumberFormat = "General"
$ws.Range("E26").Value = "  -2.27%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.396"
$c.NumberFormat = "General"
$ws.Range("E27").Value = "  -7.58%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "103.11"
$c.NumberFormat = "General"
$ws.Range("E28").Value = "  -2.69%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.677"
$c.NumberFormat = "General"
$ws.Range("E29").Value = "  -6.39%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "3.761"
$c.NumberFormat = "General"
$ws.Range("E30").Value = "  -5.08%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.07891"
$c.NumberFormat = "General"
$ws.Range("E31").Value = "  -4.39%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.596"
$c.NumberFormat = "General"
$ws.Range("E32").Value = "  -2.18%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.04572"
$c.NumberFormat = "General"
$ws.Range("E33").Value = "  -2.43%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "2.593"
$c.NumberFormat = "General"
$ws.Range("E34").Value = "  -2.45%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.9366"
$c.NumberFormat = "General"
$ws.Range("E35").Value = "  -6.53%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.5791"
$c.NumberFormat = "General"
$ws.Range("E36").Value = "  -7.55%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.595"
$c.NumberFormat = "General"
$ws.Range("E37").Value = "  -5.01%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.01563"
$c.NumberFormat = "General"
$ws.Range("E38").Value = "  -3.24%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.8452"
$c.NumberFormat = "General"
$ws.Range("E39").Value = "  +10.84%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.9995"
$c.NumberFormat = "General"
$ws.Range("E40").Value = "  +0.01%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.829"
$c.NumberFormat = "General"
$ws.Range("E41").Value = "  -5.22%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "98.35"
$c.NumberFormat = "General"
$ws.Range("E42").Value = "  -2.40%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.3730"
$c.NumberFormat = "General"
$ws.Range("E43").Value = "  -4.07%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "4.783"
$c.NumberFormat = "General"
$ws.Range("E44").Value = "  -4.82%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.1143"
$c.NumberFormat = "General"
$ws.Range("E45").Value = "  +0.80%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "6.151"
$c.NumberFormat = "General"
$ws.Range("E46").Value = "  -3.48%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.05191"
$c.NumberFormat = "General"
$ws.Range("E47").Value = "  -0.77%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "29.83"
$c.NumberFormat = "General"
$ws.Range("E48").Value = "  -3.23%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.NumberFormat = "General"
$ws.Range("E49").Value = "  -0.11%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "50.49"
$c.NumberFormat = "General"
$ws.Range("E50").Value = "  -9.30%  "
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.3337"
$c.NumberFormat = "General"
$ws.Range("E51").Value = "  -3.25%  "
